# Re-label each subgroup-level row in column A with a two-space indent so
# that the sub-rows visually nest under their category header rows
# (Sex, Age, Body-mass index, Race, Baseline Statin Treatment,
# Intensity of statin treatment, Metabolic disease, Renal function).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value  = "  Male"
$ws.Range("A5").Value  = "  Female"
$ws.Range("A7").Value  = "  <65 yr"
$ws.Range("A8").Value  = "  >65 yr"
$ws.Range("A9").Value  = "  <75 yr"
$ws.Range("A10").Value = "  >75 yr"
$ws.Range("A12").Value = "  <Median"
$ws.Range("A13").Value = "  >Median"
$ws.Range("A15").Value = "  White"
$ws.Range("A16").Value = "  Black"
$ws.Range("A17").Value = "  Other"
$ws.Range("A19").Value = "  Yes"
$ws.Range("A20").Value = "  No"
$ws.Range("A22").Value = "  High"
$ws.Range("A23").Value = "  Not High"
$ws.Range("A25").Value = "  Diabetes"
$ws.Range("A26").Value = "  Metabolic syndrome"
$ws.Range("A27").Value = "  Neither"
$ws.Range("A29").Value = "  Normal"
$ws.Range("A30").Value = "  Mild Impairment"
$ws.Range("A31").Value = "  Moderate impairment"

# Move the active selection to A31 (matches the saved view state).
$ws.Range("A31").Select()
